$wb = $excel.ActiveWorkbook

# The workbook keeps the same "想去人数" (column F) data duplicated on both
# the "展览" sheet and the "全部类型" sheet. Update the figures on both.
$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    2 = 1897
    3 = 359
    4 = 1165
    5 = 1212
    7 = 6000
    8 = 101
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
